$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 13.41179440177212
$ws.Cells.Item(3, 6).Value = 13.28345529631224
$ws.Cells.Item(4, 6).Value = 8.23763560135623
$ws.Cells.Item(5, 3).Value = 22
$ws.Cells.Item(5, 4).Value = "Khushi"
$ws.Cells.Item(5, 6).Value = 8.166612813012238
$ws.Cells.Item(5, 7).Value = "Asian"
$ws.Cells.Item(6, 3).Value = 21
$ws.Cells.Item(6, 4).Value = "Bri"
$ws.Cells.Item(6, 6).Value = 8.005597717419404
$ws.Cells.Item(6, 7).Value = "Black or African American"
$ws.Cells.Item(7, 3).Value = 30
$ws.Cells.Item(7, 4).Value = "Shadaisia"
$ws.Cells.Item(7, 6).Value = 5.361927025870195
$ws.Cells.Item(8, 3).Value = 32
$ws.Cells.Item(8, 4).Value = "Kellie"
$ws.Cells.Item(8, 6).Value = 5.237683582040133
$ws.Cells.Item(8, 7).Value = "White"
$ws.Cells.Item(9, 3).Value = 33
$ws.Cells.Item(9, 4).Value = "Shaniek"
$ws.Cells.Item(9, 6).Value = 5.010962683506764
$ws.Cells.Item(9, 7).Value = "Black or African American"
$ws.Cells.Item(10, 3).Value = 35
$ws.Cells.Item(10, 4).Value = "Lori"
$ws.Cells.Item(10, 6).Value = 4.477479880056773
$ws.Cells.Item(11, 3).Value = 34
$ws.Cells.Item(11, 4).Value = "Tina"
$ws.Cells.Item(11, 6).Value = 4.389849861394186
$ws.Cells.Item(12, 6).Value = 2.382527936458554
$ws.Cells.Item(13, 6).Value = 1.029173221199296
$ws.Cells.Item(14, 6).Value = 14.2745467615059
$ws.Cells.Item(15, 6).Value = 13.07194186949775
$ws.Cells.Item(16, 6).Value = 8.189927172263737
$ws.Cells.Item(17, 6).Value = 7.39607034879652
$ws.Cells.Item(18, 6).Value = 6.323612713011084
$ws.Cells.Item(19, 6).Value = 6.252130279629233
$ws.Cells.Item(20, 3).Value = 33
$ws.Cells.Item(20, 4).Value = "Brennan"
$ws.Cells.Item(20, 6).Value = 5.331365905335693
$ws.Cells.Item(21, 3).Value = 32
$ws.Cells.Item(21, 4).Value = "Jamarii"
$ws.Cells.Item(21, 6).Value = 5.299930968965304
$ws.Cells.Item(21, 7).Value = "Black or African American"
$ws.Cells.Item(22, 3).Value = 30
$ws.Cells.Item(22, 4).Value = "Matthew"
$ws.Cells.Item(22, 6).Value = 5.024612660958182
$ws.Cells.Item(22, 7).Value = "White"
$ws.Cells.Item(23, 6).Value = 3.207723512647401
$ws.Cells.Item(24, 6).Value = 1.345698146888841
$ws.Cells.Item(25, 6).Value = 0.1471275131748038
